$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 0.08115769241758861
$ws.Range("H2").Value = -3.231332856507831
$ws.Range("I2").Value = -16.72854336276172
$ws.Range("G3").Value = 0.117034059137713
$ws.Range("H3").Value = 0.5829870206766905
$ws.Range("G4").Value = -0.02036349073736977
$ws.Range("H4").Value = -32.01376160450518
$ws.Range("G5").Value = -0.07276019025014524
$ws.Range("H5").Value = -1.264432918329105
$ws.Range("G6").Value = 0.05313627488391683
$ws.Range("H6").Value = 49.48735605300506
$ws.Range("G7").Value = 0.02226126798619579
$ws.Range("H7").Value = 9.687904118253551
$ws.Range("G8").Value = -0.1548435155927191
$ws.Range("H8").Value = -9.741902016473572
$ws.Range("G9").Value = -0.1355902049952745
$ws.Range("H9").Value = 0.9473544793316588
$ws.Range("G10").Value = -0.06593626870677871
$ws.Range("H10").Value = 38.9632281303391
$ws.Range("G11").Value = -0.06650545617957808
$ws.Range("H11").Value = 0.05320859499277958
$ws.Range("G12").Value = -0.3161817590787152
$ws.Range("H12").Value = 23.84764959420039
$ws.Range("G13").Value = -0.3663680869026293
$ws.Range("H13").Value = 18.385255879795
$ws.Range("G14").Value = -0.04651163650635912
$ws.Range("H14").Value = 8.325616555093289
$ws.Range("G15").Value = 0.03494418024063498
$ws.Range("H15").Value = 142.2585231748803
$ws.Range("G16").Value = 0.1003747225743274
$ws.Range("H16").Value = -31.26087085192736
$ws.Range("G17").Value = 0.1838908558506729
$ws.Range("H17").Value = 49.96626033740213
$ws.Range("G18").Value = 0.1296040433581075
$ws.Range("H18").Value = -6.394133727838571
$ws.Range("G19").Value = 0.1287065971766345
$ws.Range("H19").Value = 34.90303797102566
$ws.Range("G20").Value = 0.03553634587461138
$ws.Range("H20").Value = 38.62945145938912
$ws.Range("G21").Value = 0.07429486866570362
$ws.Range("H21").Value = -1.160406235556423
$ws.Range("G24").Value = 0.07360093579903874
$ws.Range("H24").Value = -26.71894113235834
$ws.Range("G25").Value = 0.151198295228831
$ws.Range("H25").Value = -0.2297899960804376
$ws.Range("G26").Value = 0.07258594491379972
$ws.Range("H26").Value = -8.263874144986374
$ws.Range("G27").Value = 0.08625796586418241
$ws.Range("H27").Value = -13.61472044110182
$ws.Range("G28").Value = -0.239281318743601
$ws.Range("H28").Value = -12.27582313651761
$ws.Range("G29").Value = -0.2308367661083208
$ws.Range("H29").Value = -12.43509990849802
$ws.Range("G30").Value = 0.06680726225209016
$ws.Range("H30").Value = 51.38177106050452
$ws.Range("G31").Value = 0.01369019412764048
$ws.Range("H31").Value = -48.01148377666485
$ws.Range("G32").Value = 0.09331645183010881
$ws.Range("H32").Value = -1.731305945303205
$ws.Range("G33").Value = 0.1379215833466158
$ws.Range("H33").Value = 32.66169702608172
$ws.Range("G34").Value = 0.04742064126246821
$ws.Range("H34").Value = 2.138701995910228
$ws.Range("G35").Value = 0.001307587305217545
$ws.Range("H35").Value = -82.74160164419173
$ws.Range("G36").Value = 0.04654054422263969
$ws.Range("H36").Value = -19.39111476020775
$ws.Range("G37").Value = 0.0721150862937399
$ws.Range("H37").Value = 2.543757457111079
$ws.Range("G38").Value = 0.02837899317752992
$ws.Range("H38").Value = -45.82756780583328
$ws.Range("G39").Value = 0.0436362556584001
$ws.Range("H39").Value = 110.4237000443578
$ws.Range("G40").Value = -0.001190059719832217
$ws.Range("H40").Value = 85.98444013014584
$ws.Range("G41").Value = 0.01253345651831934
$ws.Range("H41").Value = -64.55054793924863
$ws.Range("G42").Value = 0.1490893523402751
$ws.Range("H42").Value = 11.52193353248969
$ws.Range("G43").Value = 0.1496959682699459
$ws.Range("H43").Value = 0.4830524134918599
$ws.Range("G44").Value = -0.003716466556209904
$ws.Range("H44").Value = 56.33584287181515
$ws.Range("G45").Value = 0.01396031570668361
$ws.Range("H45").Value = 227.1446928598841
$ws.Range("G46").Value = -0.002975050573741749
$ws.Range("H46").Value = 9.65657156881816
$ws.Range("G47").Value = -0.005447829636818419
$ws.Range("H47").Value = 41.28532912987693
$ws.Range("G48").Value = 0.05920818471777103
$ws.Range("H48").Value = 17.775407373198
$ws.Range("G49").Value = 0.06499854662549265
$ws.Range("H49").Value = -1.613926435931025
$ws.Range("G50").Value = 0.1349993904718859
$ws.Range("H50").Value = -16.28343517921649
$ws.Range("G51").Value = 0.1629246195032314
$ws.Range("H51").Value = -4.786454669877129
$ws.Range("G52").Value = -0.1796242364698385
$ws.Range("H52").Value = -11.96354486289387
$ws.Range("G53").Value = -0.1521467331565825
$ws.Range("H53").Value = -20.69916800198779
$ws.Range("G54").Value = 0.09887457216937237
$ws.Range("H54").Value = 5.497274113789521
$ws.Range("G55").Value = 0.1169786815214399
$ws.Range("H55").Value = 3.447036736077018
$ws.Range("G56").Value = -0.02364835448933883
$ws.Range("H56").Value = -223.9158282821367
$ws.Range("G57").Value = -0.0159995907739671
$ws.Range("H57").Value = 30.03316117168243
$ws.Range("G58").Value = 0.03886477168281564
$ws.Range("H58").Value = -31.07613216710932
$ws.Range("G59").Value = 0.05891662056366272
$ws.Range("H59").Value = -17.9661914603764
$ws.Range("G60").Value = 0.05481226376046851
$ws.Range("H60").Value = -21.66761907482947
$ws.Range("G61").Value = 0.07507534644010314
$ws.Range("H61").Value = 57.96428393497387
$ws.Range("G62").Value = 0.06589821670159328
$ws.Range("H62").Value = -9.679384743854531
$ws.Range("G63").Value = 0.07190248407717978
$ws.Range("H63").Value = 9.953556672252937
$ws.Range("G64").Value = -0.007012464004576278
$ws.Range("H64").Value = 83.0680235053599
$ws.Range("G65").Value = -0.01909697371153306
$ws.Range("H65").Value = 61.28236378456542
$ws.Range("G66").Value = 0.04328298747059912
$ws.Range("H66").Value = 128.5995758385288
$ws.Range("G67").Value = 0.04010988899251833
$ws.Range("H67").Value = 53.38456624910132
$ws.Range("G68").Value = -0.0154335170111205
$ws.Range("H68").Value = -2807.856055275632
$ws.Range("G69").Value = 0.005293215960864637
$ws.Range("H69").Value = 140.9935733864504
$ws.Range("G70").Value = -0.03554354903007304
$ws.Range("H70").Value = -29.49380946922368
$ws.Range("G71").Value = -0.0597394838831968
$ws.Range("H71").Value = -8.435956269299934
$ws.Range("G72").Value = -0.155372895282417
$ws.Range("H72").Value = -4.745198112503925
$ws.Range("G73").Value = -0.1540640149776858
$ws.Range("H73").Value = -6.405269341099813
$ws.Range("G74").Value = 0.149325749823856
$ws.Range("H74").Value = 18.46445064533247
$ws.Range("G75").Value = 0.1403727918048235
$ws.Range("H75").Value = 3.847553832025317
$ws.Range("G76").Value = -0.01962171400663281
$ws.Range("H76").Value = 43.02298240146806
$ws.Range("G77").Value = -0.02952884765206244
$ws.Range("H77").Value = 36.06676631979698
$ws.Range("G78").Value = 0.1055010152631173
$ws.Range("H78").Value = 14.46159774493277
$ws.Range("G79").Value = 0.105348857773641
$ws.Range("H79").Value = 9.164789999100844
$ws.Range("G80").Value = -0.202407080733238
$ws.Range("H80").Value = -24.62503445642676
$ws.Range("G81").Value = -0.2163432841009074
$ws.Range("H81").Value = 0.04859074454615082
$ws.Range("G82").Value = 0.1665639310621269
$ws.Range("H82").Value = 20.04929585060155
$ws.Range("G83").Value = 0.1812627777629491
$ws.Range("H83").Value = 10.11008609655441
$ws.Range("G84").Value = 0.04560335558397417
$ws.Range("H84").Value = 225.8539891274808
$ws.Range("G85").Value = 0.08368413833058931
$ws.Range("H85").Value = 269.6429139113765
